# "Changes in Create Call" -- expand the CreateCallData sheet with the full
# set of customer/call columns (C..S), refresh the ExpenseData amount, and
# switch the active tab from ExpenseData to CreateCallData.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. CreateCallData sheet: rewrite the header row + data row with the new,
#    much wider set of columns (A..S).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CreateCallData")

$headers = @("Customer First Name", "City ", "Customer Middle Name", "Customer Last Name", `
    "Address 1", "Address 2", "State", "Customer Id", "Primary Number", "Email Id", `
    "Org Name", "Sub Org", "Customer Type", "Problem Description", "Product Name", `
    "ProductSerialNo", "WarrantyType", "SepcialInstructions", "Remarks")

$data = @("Steve", "Bangalore", "Martin", "Ontoyo", `
    "Los Angeles Buckking Ham Street", "New York USA", "New York", "cust883", "ph-932883832", `
    "stevel@jobs.co.in", "Apple", "Ottovia", "Indiviual", "Remarks is there", "WildCraft", `
    "Craft920393", "2  years", "Warranty Type Bonues for  1 year", "Do renewal after 3 years")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $data[$i]
}

# Whole used range gets a Text number format.
$ws.Range("A1:S2").NumberFormat = "@"

# Column widths: A-D and F-M are 23 chars wide, E (longest header) auto-fits
# wider to fit "Los Angeles Buckking Ham Street".
$ws.Range("A1:D2").ColumnWidth = 22.16
$ws.Range("F1:M2").ColumnWidth = 22.16
$ws.Columns.Item(5).ColumnWidth = 29.88

# J2 (the email address) becomes a mailto hyperlink.
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:stevel@jobs.co.in")

# CreateCallData becomes the active sheet / tab, with C16 selected.
$ws.Activate() | Out-Null
$ws.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. ExpenseData sheet: bump the amount from Rs400 to Rs600. (It stops
#    being the active tab because CreateCallData now is.)
# ---------------------------------------------------------------------
$wsExpense = $wb.Worksheets.Item("ExpenseData")
$wsExpense.Range("A2").Value = "Rs600"
